$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 17 ("Programa:" -> becomes "Avaliação:") loses its B/C paragraph text.
$ws.Range("B17:C17").Clear()
$ws.Range("A17").Value = "Avaliação:"

# Row 14 ("7459752..." data row -> becomes "Short syllabus:" label) loses B/C.
$ws.Range("B14:C14").Clear()
$ws.Range("A14").Value = "Short syllabus:"

# Row 13 gains a label in column A ("Programa resumido:"); B/C keep holding
# the exact text already used by B8/C8 ("01/01/2022") - copy those cells so
# the date-looking text is stored as the same shared string/style, not
# re-parsed as a date.
$ws.Range("B8").Copy($ws.Range("B13"))
$ws.Range("C8").Copy($ws.Range("C13"))
$ws.Range("A13").Value = "Programa resumido:"

# Row 15 keeps its A/B/C cells, just changes text.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C15").Value = "984972 - Hugo Ricardo Zschommler Sandim"

# Row 16 label only.
$ws.Range("A16").Value = "Syllabus:"

# Row 10 keeps A/B/C, B/C text changes.
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C10").Value = "984972 - Hugo Ricardo Zschommler Sandim"

# Row 18 ("Syllabus:" label-only row) gains B/C cells. Seed them via Copy
# from an existing, correctly styled column B/C cell so the new cells pick
# up the right style (s=2 / s=3) instead of defaulting incorrectly, then
# overwrite with the real text.
$ws.Range("B10").Copy($ws.Range("B18"))
$ws.Range("C10").Copy($ws.Range("C18"))
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C18").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

# Row 19 ("Avaliação:" label-only row) gains B/C cells the same way.
$ws.Range("B10").Copy($ws.Range("B19"))
$ws.Range("C10").Copy($ws.Range("C19"))
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras"
$ws.Range("C19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras"

# Row 20 keeps its A/B/C cells, just changes text.
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média Aritmética dos Projetos, Trabalhos, Relatórios e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
$ws.Range("C20").Value = "Média Aritmética dos Projetos, Trabalhos, Relatórios e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."

# Row 21 keeps its A/B/C cells, just changes text.
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."

# Drop the two now-redundant trailing rows (their content has been folded
# into rows 13-21 above). Delete bottom-up so row numbers stay valid.
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()

# ---------------------------------------------------------------------------
# Row heights
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).EntireRow.AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
